# Updates the cryptos worksheet's Price (column D) and Volume(1h) (column E)
# values for rows 2-51, matching the latest scrape snapshot.
#
# Column D values are plain numeric-looking strings (e.g. "109.68") that
# Excel would otherwise silently convert to numbers on assignment, but the
# workbook stores them as literal text. To force text storage we briefly mark
# the destination cell's number format as Text ("@") before assigning the
# value, then restore the cell's original style so no visible formatting
# change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue($cell, [string]$value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$updates = @(
    @{ Row = 2; D = '42.373.05'; E = '  -3.02%  ' },
    @{ Row = 3; D = '2.219.63'; E = '  -2.36%  ' },
    @{ Row = 4; D = '1.00'; E = '  +0.09%  ' },
    @{ Row = 5; D = '109.68'; E = '  -7.30%  ' },
    @{ Row = 6; D = '286.79'; E = $null },
    @{ Row = 7; D = '0.619'; E = '  -3.49%  ' },
    @{ Row = 8; D = $null; E = '  -0.31%  ' },
    @{ Row = 9; D = '0.594'; E = '  -4.50%  ' },
    @{ Row = 10; D = '43.14'; E = $null },
    @{ Row = 11; D = '0.0904'; E = '  -4.25%  ' },
    @{ Row = 12; D = '54.19'; E = '  +0.26%  ' },
    @{ Row = 13; D = '8.56'; E = '  -8.99%  ' },
    @{ Row = 14; D = $null; E = '  +11.53%  ' },
    @{ Row = 15; D = $null; E = '  -3.22%  ' },
    @{ Row = 16; D = '14.78'; E = '  -5.96%  ' },
    @{ Row = 17; D = '2.550.23'; E = $null },
    @{ Row = 18; D = '2.227.57'; E = '  -1.90%  ' },
    @{ Row = 19; D = '42.251.77'; E = '  -3.21%  ' },
    @{ Row = 20; D = '7.12'; E = '  +2.96%  ' },
    @{ Row = 21; D = $null; E = '  -5.41%  ' },
    @{ Row = 22; D = '72.68'; E = '  +0.36%  ' },
    @{ Row = 23; D = '3.34'; E = '  +12.74%  ' },
    @{ Row = 24; D = '2.37'; E = '  -1.18%  ' },
    @{ Row = 25; D = '228.68'; E = '  -2.49%  ' },
    @{ Row = 26; D = '8.86'; E = '  -8.82%  ' },
    @{ Row = 27; D = $null; E = '  -1.83%  ' },
    @{ Row = 28; D = '11.36'; E = '  -7.35%  ' },
    @{ Row = 29; D = $null; E = '  -2.68%  ' },
    @{ Row = 30; D = $null; E = '  -4.81%  ' },
    @{ Row = 31; D = '172.38'; E = '  -1.13%  ' },
    @{ Row = 32; D = '36.62'; E = '  -12.52%  ' },
    @{ Row = 33; D = '20.71'; E = '  -3.80%  ' },
    @{ Row = 34; D = '0.0868'; E = '  -5.64%  ' },
    @{ Row = 35; D = '5.53'; E = '  -3.74%  ' },
    @{ Row = 36; D = '4.94'; E = '  +7.17%  ' },
    @{ Row = 37; D = '0.125'; E = '  -4.10%  ' },
    @{ Row = 38; D = '4.14'; E = '  -3.21%  ' },
    @{ Row = 39; D = '0.0367'; E = '  -4.57%  ' },
    @{ Row = 40; D = '0.104'; E = '  -4.66%  ' },
    @{ Row = 41; D = '73.71'; E = '  +2.48%  ' },
    @{ Row = 42; D = '2.37'; E = '  -6.96%  ' },
    @{ Row = 43; D = '0.228'; E = '  -4.99%  ' },
    @{ Row = 44; D = $null; E = '  -0.02%  ' },
    @{ Row = 45; D = '12.17'; E = '  -12.48%  ' },
    @{ Row = 46; D = $null; E = '  -7.02%  ' },
    @{ Row = 47; D = '5.34'; E = '  -7.04%  ' },
    @{ Row = 48; D = '1.69'; E = '  +8.64%  ' },
    @{ Row = 49; D = '1.26'; E = '  -0.89%  ' },
    @{ Row = 50; D = '8.38'; E = '  -2.11%  ' },
    @{ Row = 51; D = '100.54'; E = '  -2.56%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextCellValue $ws.Cells.Item($r, 4) $u.D
    }
    if ($null -ne $u.E) {
        Set-TextCellValue $ws.Cells.Item($r, 5) $u.E
    }
}
